$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2891508418385627
$ws.Range("C2").Value2 = 0.1067110064959991
$ws.Range("D2").Value2 = 0.08690141491899439
$ws.Range("E2").Value2 = 0.1362578027949581
$ws.Range("F2").Value2 = 2.073173598615682
$ws.Range("I2").Value2 = 1.503774068279739
$ws.Range("K2").Value2 = 0.3934059789598336
$ws.Range("M2").Value2 = 0.2581166958739516

$ws.Range("B3").Value2 = 0.2669171553448564
$ws.Range("C3").Value2 = 0.09709692343255938
$ws.Range("D3").Value2 = 0.08603626890612048
$ws.Range("E3").Value2 = 0.1251582388022427
$ws.Range("F3").Value2 = 2.010269201892442
$ws.Range("I3").Value2 = 1.466216486311239
$ws.Range("K3").Value2 = 0.3615412303506957
$ws.Range("M3").Value2 = 0.237157190233944

$ws.Range("B4").Value2 = 0.2535071141166725
$ws.Range("C4").Value2 = 0.09125632574800591
$ws.Range("D4").Value2 = 0.08549156246950673
$ws.Range("E4").Value2 = 0.1184174020751314
$ws.Range("F4").Value2 = 1.972309074682485
$ws.Range("I4").Value2 = 1.443527407694177
$ws.Range("K4").Value2 = 0.3422739594769268
$ws.Range("M4").Value2 = 0.2244549096416009

$ws.Range("B5").Value2 = 0.2481029561334651
$ws.Range("C5").Value2 = 0.08889176286535871
$ws.Range("D5").Value2 = 0.08526616324851943
$ws.Range("E5").Value2 = 0.1156888395406384
$ws.Range("F5").Value2 = 1.957005840836317
$ws.Range("I5").Value2 = 1.434374161633727
$ws.Range("K5").Value2 = 0.3344968874324934
$ws.Range("M5").Value2 = 0.2193201731429184

$ws.Range("B6").Value2 = 0.2472092509667334
$ws.Range("C6").Value2 = 0.08850006149656053
$ws.Range("D6").Value2 = 0.08522852797962699
$ws.Range("E6").Value2 = 0.1152368650858122
$ws.Range("F6").Value2 = 1.954474737007118
$ws.Range("I6").Value2 = 1.432859853972857
$ws.Range("K6").Value2 = 0.333209998776141
$ws.Range("M6").Value2 = 0.2184700507953963

$ws.Range("B7").Value2 = 0.253433986867492
$ws.Range("C7").Value2 = 0.09122437379535597
$ws.Range("D7").Value2 = 0.08548853657422484
$ws.Range("E7").Value2 = 0.1183805297003602
$ws.Range("F7").Value2 = 1.972102019826593
$ws.Range("I7").Value2 = 1.443403589210021
$ws.Range("K7").Value2 = 0.3421687740370203
$ws.Range("M7").Value2 = 0.2243854931377172

$ws.Range("B8").Value2 = 0.2814344113428433
$ws.Range("C8").Value2 = 0.1033829862607831
$ws.Range("D8").Value2 = 0.08660589618328629
$ws.Range("E8").Value2 = 0.1324150529918455
$ws.Range("F8").Value2 = 2.051345845206299
$ws.Range("I8").Value2 = 1.490746713267058
$ws.Range("K8").Value2 = 0.3823569544608745
$ws.Range("M8").Value2 = 0.2508549414455175

$ws.Range("B9").Value2 = 0.3382718547084096
$ws.Range("C9").Value2 = 0.1277314873175897
$ws.Range("D9").Value2 = 0.08869135224746927
$ws.Range("E9").Value2 = 0.1605420427246003
$ws.Range("F9").Value2 = 2.212060838129958
$ws.Range("I9").Value2 = 1.586568621275859
$ws.Range("K9").Value2 = 0.4635513476984556
$ws.Range("M9").Value2 = 0.3041085383462701

$ws.Range("B10").Value2 = 0.3812282067259787
$ws.Range("C10").Value2 = 0.1459438294195934
$ws.Range("D10").Value2 = 0.09016119167802117
$ws.Range("E10").Value2 = 0.1815995221152704
$ws.Range("F10").Value2 = 2.333468205380086
$ws.Range("I10").Value2 = 1.658843116397705
$ws.Range("K10").Value2 = 0.5246972720037775
$ws.Range("M10").Value2 = 0.3440906120354654

$ws.Range("B11").Value2 = 0.4010354442656592
$ws.Range("C11").Value2 = 0.1543027926482239
$ws.Range("D11").Value2 = 0.09081677006545164
$ws.Range("E11").Value2 = 0.1912696009673454
$ws.Range("F11").Value2 = 2.38944248506354
$ws.Range("I11").Value2 = 1.692142178706519
$ws.Range("K11").Value2 = 0.5528469701608287
$ws.Range("M11").Value2 = 0.3624734781092513

$ws.Range("B12").Value2 = 0.4085745302640191
$ws.Range("C12").Value2 = 0.1574790098156598
$ws.Range("D12").Value2 = 0.09106317665607122
$ws.Range("E12").Value2 = 0.1949448676130103
$ws.Range("F12").Value2 = 2.410747017387337
$ws.Range("I12").Value2 = 1.70481306698548
$ws.Range("K12").Value2 = 0.5635551530296254
$ws.Range("M12").Value2 = 0.3694631709900165

$ws.Range("B13").Value2 = 0.4069491363540862
$ws.Range("C13").Value2 = 0.1567944683270639
$ws.Range("D13").Value2 = 0.09101019022726575
$ws.Range("E13").Value2 = 0.1941527309334745
$ws.Range("F13").Value2 = 2.406153871573025
$ws.Range("I13").Value2 = 1.702081428713882
$ws.Range("K13").Value2 = 0.5612467896604301
$ws.Range("M13").Value2 = 0.3679565401984348

$ws.Range("B14").Value2 = 0.4016549161150635
$ws.Range("C14").Value2 = 0.1545638831154008
$ws.Range("D14").Value2 = 0.09083707891758763
$ws.Range("E14").Value2 = 0.1915716964803522
$ws.Range("F14").Value2 = 2.391193045236207
$ws.Range("I14").Value2 = 1.693183388554402
$ws.Range("K14").Value2 = 0.5537269640742011
$ws.Range("M14").Value2 = 0.3630479503521542

$ws.Range("B15").Value2 = 0.3984170763596637
$ws.Range("C15").Value2 = 0.1531990056478207
$ws.Range("D15").Value2 = 0.09073080366664499
$ws.Range("E15").Value2 = 0.1899924959595225
$ws.Range("F15").Value2 = 2.382043241879217
$ws.Range("I15").Value2 = 1.687741080586875
$ws.Range("K15").Value2 = 0.5491271840840852
$ws.Range("M15").Value2 = 0.3600450234657302

$ws.Range("B16").Value2 = 0.3799391336107476
$ws.Range("C16").Value2 = 0.1453990597689199
$ws.Range("D16").Value2 = 0.09011808797378507
$ws.Range("E16").Value2 = 0.1809694185856117
$ws.Range("F16").Value2 = 2.329825248187319
$ws.Range("I16").Value2 = 1.656675477386472
$ws.Range("K16").Value2 = 0.5228643851457946
$ws.Range("M16").Value2 = 0.3428932101401045

$ws.Range("B17").Value2 = 0.3686718564833882
$ws.Range("C17").Value2 = 0.1406331649766059
$ws.Range("D17").Value2 = 0.08973888422915621
$ws.Range("E17").Value2 = 0.1754575837911574
$ws.Range("F17").Value2 = 2.29798290400521
$ws.Range("I17").Value2 = 1.637726077835751
$ws.Range("K17").Value2 = 0.5068389054073634
$ws.Range("M17").Value2 = 0.332421369719377

$ws.Range("B18").Value2 = 0.3622162558482671
$ws.Range("C18").Value2 = 0.1378989083177373
$ws.Range("D18").Value2 = 0.08951954523979211
$ws.Range("E18").Value2 = 0.1722958598637447
$ws.Range("F18").Value2 = 2.279738130220437
$ws.Range("I18").Value2 = 1.626866471111526
$ws.Range("K18").Value2 = 0.4976528933655686
$ws.Range("M18").Value2 = 0.3264165520435967

$ws.Range("B19").Value2 = 0.3600347937291133
$ws.Range("C19").Value2 = 0.1369743251179329
$ws.Range("D19").Value2 = 0.08944506844676425
$ws.Range("E19").Value2 = 0.1712268104412189
$ws.Range("F19").Value2 = 2.273572770004961
$ws.Range("I19").Value2 = 1.623196372285889
$ws.Range("K19").Value2 = 0.4945480501415034
$ws.Range("M19").Value2 = 0.3243865536811441

$ws.Range("B20").Value2 = 0.369868684202288
$ws.Range("C20").Value2 = 0.1411397807903541
$ws.Range("D20").Value2 = 0.08977937829553895
$ws.Range("E20").Value2 = 0.1760434419233832
$ws.Range("F20").Value2 = 2.301365312926407
$ws.Range("I20").Value2 = 1.639739171565282
$ws.Range("K20").Value2 = 0.5085415899376358
$ws.Range("M20").Value2 = 0.333534216702617

$ws.Range("B21").Value2 = 0.4032089097061373
$ws.Range("C21").Value2 = 0.1552187636825977
$ws.Range("D21").Value2 = 0.09088797577648933
$ws.Range("E21").Value2 = 0.1923294419643895
$ws.Range("F21").Value2 = 2.395584454954502
$ws.Range("I21").Value2 = 1.695795290471281
$ws.Range("K21").Value2 = 0.555934399069173
$ws.Range("M21").Value2 = 0.3644889447789339

$ws.Range("B22").Value2 = 0.4252232854565818
$ws.Range("C22").Value2 = 0.1644835867646748
$ws.Range("D22").Value2 = 0.09160175861713071
$ws.Range("E22").Value2 = 0.2030516544856837
$ws.Range("F22").Value2 = 2.457793805870011
$ws.Range("I22").Value2 = 1.732788585509965
$ws.Range("K22").Value2 = 0.5871912726125572
$ws.Range("M22").Value2 = 0.3848860016509192

$ws.Range("B23").Value2 = 0.4134531660330083
$ws.Range("C23").Value2 = 0.1595329066427951
$ws.Range("D23").Value2 = 0.09122177261780706
$ws.Range("E23").Value2 = 0.197321722826203
$ws.Range("F23").Value2 = 2.424533341383324
$ws.Range("I23").Value2 = 1.713011630923475
$ws.Range("K23").Value2 = 0.5704828484828113
$ws.Range("M23").Value2 = 0.3739843250614001

$ws.Range("B24").Value2 = 0.3693275295492242
$ws.Range("C24").Value2 = 0.1409107219810721
$ws.Range("D24").Value2 = 0.08976107507302089
$ws.Range("E24").Value2 = 0.1757785533652694
$ws.Range("F24").Value2 = 2.299835933385992
$ws.Range("I24").Value2 = 1.63882894399066
$ws.Range("K24").Value2 = 0.5077717214200845
$ws.Range("M24").Value2 = 0.3330310500839602

$ws.Range("B25").Value2 = 0.3226870816314147
$ws.Range("C25").Value2 = 0.1210888805490242
$ws.Range("D25").Value2 = 0.08813827952818443
$ws.Range("E25").Value2 = 0.1528657194833087
$ws.Range("F25").Value2 = 2.168004768815962
$ws.Range("I25").Value2 = 1.560321205311254
$ws.Range("K25").Value2 = 0.4413268219244912
$ws.Range("M25").Value2 = 0.2895541743523822

